# Update the bilingual/trilingual title in row 1: drop the period after
# "6.4.1.2" in the Russian and English titles (Kyrgyz title text is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "6.4.1.2 Потери воды при транспортировке"
$ws.Range("C1").Value = "6.4.1.2 Percentage of water loss during transportation"

# Updated 2022 data points.
$ws.Range("P5").Value = 2388
$ws.Range("P10").Value = 335.3
$ws.Range("P16").Value = 27.3
$ws.Range("P21").Value = 24.3

# Move the active selection to S3.
$ws.Range("S3").Select()
